# Add a new "staff_information" worksheet between "customer_details" and "menu",
# populate it with staff data, and clear the now-unused cell formatting that used
# to live on the customer_details sheet (columns D:G).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet right after "customer_details" -------------------
$afterSheet = $wb.Worksheets.Item("customer_details")
$staff = $wb.Worksheets.Add($null, $afterSheet)
$staff.Name = "staff_information"

# --- 2. Header row --------------------------------------------------------------
$staff.Cells.Item(1, 1).Value = "staff_id"
$staff.Cells.Item(1, 2).Value = "role"
$staff.Cells.Item(1, 3).Value = "salary"
$staff.Cells.Item(1, 4).Value = "hire_date"

# --- 3. Data rows -----------------------------------------------------------------
$rows = @(
    @(1,  "Waiter",        2000, 44562),
    @(2,  "Chef",          3000, 44576),
    @(3,  "Manager",       4000, 44593),
    @(4,  "Bartender",     2500, 44607),
    @(5,  "Host/Hostess",  1800, 44621),
    @(6,  "Sous Chef",     3200, 44635),
    @(7,  "Server",        2000, 44652),
    @(8,  "Dishwasher",    1800, 44666),
    @(9,  "Line Cook",     2400, 44682),
    @(10, "Busser",        1700, 44696),
    @(11, "Food Runner",   1700, 44713),
    @(12, "Prep Cook",     2000, 44727),
    @(13, "Head Chef",     4000, 44743),
    @(14, "Floor Manager", 3500, 44757),
    @(15, "Sommelier",     3000, 44774)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $staff.Cells.Item($r, 1).Value = $row[0]
    $staff.Cells.Item($r, 2).Value = $row[1]
    $staff.Cells.Item($r, 3).Value = $row[2]
    $staff.Cells.Item($r, 4).Value = $row[3]
}

# --- 4. Give hire_date the same date format already used elsewhere in the book --
# (copy format from an existing date cell so we reuse the existing style instead
# of minting a brand new number format)
$dateSource = $wb.Worksheets.Item("bookings").Range("B2")
$dateSource.Copy()
$staff.Range("D2:D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$staff.Columns.Item(4).ColumnWidth = 9.6

# --- 5. Make the new sheet the active tab (matches the saved workbook view) -----
$staff.Activate()
$staff.Range("A1").Select()

# --- 6. Clear the now-unused custom formatting on customer_details D1:G16 -------
$customers = $wb.Worksheets.Item("customer_details")
$customers.Range("D1:G16").ClearFormats()
